# "Fixed up Myxicola and a few more loose ends in second review"
#
# On the "Materials" sheet:
#   - Drop the Taxon_Local_ID / iNaturalistTaxonId column.
#   - Drop the suborder / infraorder / superfamily columns (the Taxa sheet
#     already carries Suborder/Infraorder/Superfamily; these lowercase
#     Darwin-Core-style template columns on Materials were extraneous).
#   - Fix the "${summary.Author}" template placeholder to the correct
#     "${summary.authority}" field name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")
$headerRow = $ws.Rows.Item(1)

# Locate the columns to remove by their header text so the script doesn't
# depend on hard-coded column letters.
$colTaxonLocalId = $headerRow.Find("Taxon_Local_ID").Column
$colSuborder     = $headerRow.Find("suborder").Column
$colInfraorder   = $headerRow.Find("infraorder").Column
$colSuperfamily  = $headerRow.Find("superfamily").Column

# Delete from right to left so the remaining column numbers stay valid.
$cols = @($colTaxonLocalId, $colSuborder, $colInfraorder, $colSuperfamily) | Sort-Object -Descending
foreach ($col in $cols) {
    $ws.Columns.Item($col).Delete()
}

# Correct the mis-named template placeholder left over in row 2.
$badCell = $ws.Rows.Item(2).Find('${summary.Author}')
$badCell.Value = '${summary.authority}'
